$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the cells that should be removed from row 2 (X2, Y2, Z2, AA2, AB2, AD2)
# AC2 (value 1) remains untouched.
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AD2").ClearContents()
